$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A="<are>";      B="<arl>";      C=15}
    @{Row=3;  A="<his>";      B="<his>";      C=17}
    @{Row=4;  A="<kilo>";     B="<kilo>";     C=16}
    @{Row=5;  A="<come>";     B="<come>";     C=17}
    @{Row=6;  A="<sentence>"; B="<seeven>";   C=13}
    @{Row=7;  A="<zero>";     B="<zero>";     C=10}
    @{Row=8;  A="<shift>";    B="<shift>";    C=12}
    @{Row=9;  A="<so>";       B="<so>";       C=9}
    @{Row=10; A="<lima>";     B="<lima>";     C=10}
    @{Row=11; A="<be>";       B="<be>";       C=13}
    @{Row=12; A="<him>";      B="<hin>";      C=11}
    @{Row=13; A="<your>";     B="<your>";     C=15}
    @{Row=14; A="<are>";      B="<are>";      C=11}
    @{Row=15; A="<in>";       B="<in>";       C=10}
    @{Row=16; A="<number>";   B="<number>";   C=14}
    @{Row=17; A="<left>";     B="<that>";     C=15}
    @{Row=18; A="<no>";       B="<no>";       C=11}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
